$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134, pushing existing rows 134:208 down to 135:209
$ws.Rows("134:134").Insert()

# Populate the newly inserted row 134 with its values
$ws.Cells.Item(134, 1).Value = 5
$ws.Cells.Item(134, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(134, 3).Value = "Maule"
$ws.Cells.Item(134, 4).Value = 44518
$ws.Cells.Item(134, 5).Value = 7
$ws.Cells.Item(134, 6).Value = 100114014
$ws.Cells.Item(134, 7).Value = "Betarraga"
$ws.Cells.Item(134, 8).Value = "Sin especificar"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 6000
$ws.Cells.Item(134, 11).Value = 500
$ws.Cells.Item(134, 12).Value = 500
$ws.Cells.Item(134, 13).Value = 500
$ws.Cells.Item(134, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(134, 15).Value = "Región del Maule"
$ws.Cells.Item(134, 16).Value = 100
$ws.Cells.Item(134, 17).Value = 5
$ws.Cells.Item(134, 18).Value = "Hortaliza"
